# Regenerate orders with updated distance/size codes.
# Applies the following text substitutions to every string-valued cell
# on the active sheet: D64->D69, D80->D86, D51->D55, S30->S31

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = $ws.UsedRange.Rows.Count
$cols = $ws.UsedRange.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -is [string]) {
            $newVal = $v.Replace("D64", "D69").Replace("D80", "D86").Replace("D51", "D55").Replace("S30", "S31")
            if ($newVal -ne $v) {
                $cell.Value = $newVal
            }
        }
    }
}
